$wb = $excel.ActiveWorkbook

# Update the "Status" shared text from "Ready for handoff" to "Handback transform failed".
# This string is shared across the Overview, zh-cn and de-de sheets (row 3 / the
# 980bab84-... entry), so updating it once on each sheet keeps them all in sync.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file-name mismatch error detail for the second
# (980bab84-...) row on both locale sheets, in column K ("Error Detail").
$zhcn.Range("K3").Value = "Handback file name: bfc20kme.51l is different with handoff file name: 980bab84-0dfd-40db-b4c7-2170b0f9008a.beecac91a2fd5bcb7387141c9d7ef2719430c91a.zh-cn."

$dede.Range("K3").Value = "Handback file name: bfc20kme.51l is different with handoff file name: 980bab84-0dfd-40db-b4c7-2170b0f9008a.beecac91a2fd5bcb7387141c9d7ef2719430c91a.de-de."
